$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A70").Value = "Can I delete a curve from GEO?"
$ws.Range("B70").Value = "llama3.2:latest"
$ws.Range("C70").Value = "No, you can only edit the data values for all other curve types. You cannot delete a curve from GEO."

$ws.Range("A71").Value = "Can I delete a curve from GEO?"
$ws.Range("B71").Value = "llama3.2:latest"
$ws.Range("C71").Value = "No, you can only edit the data values for all other curve types. You cannot delete a curve from GEO."
